# Weekly update: insert two new price rows (date 44918) for "Choclo" /
# "Dulce o Americano" in "Región de O'Higgins" (Primera & Segunda),
# right before the current row 450, shifting all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 450-451 (existing rows 450+ shift down to 452+)
$ws.Rows("450:451").Insert()

# --- New row 450: Choclo / Dulce o Americano / Primera ---
$ws.Range("A450").Value = 10
$ws.Range("B450").Value = "Vega Modelo de Temuco"
$ws.Range("C450").Value = "La Araucanía"
$ws.Range("D450").Value = 44918
$ws.Range("E450").Value = 9
$ws.Range("F450").Value = 100112024
$ws.Range("G450").Value = "Choclo"
$ws.Range("H450").Value = "Dulce o Americano"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 3250
$ws.Range("K450").Value = 250
$ws.Range("L450").Value = 300
$ws.Range("M450").Value = 281
$ws.Range("N450").Value = "$/unidad"
$ws.Range("O450").Value = "Región de O'Higgins"
$ws.Range("P450").Value = 281
$ws.Range("Q450").Value = 1
$ws.Range("R450").Value = "Hortaliza"

# --- New row 451: Choclo / Dulce o Americano / Segunda ---
$ws.Range("A451").Value = 10
$ws.Range("B451").Value = "Vega Modelo de Temuco"
$ws.Range("C451").Value = "La Araucanía"
$ws.Range("D451").Value = 44918
$ws.Range("E451").Value = 9
$ws.Range("F451").Value = 100112024
$ws.Range("G451").Value = "Choclo"
$ws.Range("H451").Value = "Dulce o Americano"
$ws.Range("I451").Value = "Segunda"
$ws.Range("J451").Value = 110
$ws.Range("K451").Value = 200
$ws.Range("L451").Value = 200
$ws.Range("M451").Value = 200
$ws.Range("N451").Value = "$/unidad"
$ws.Range("O451").Value = "Región de O'Higgins"
$ws.Range("P451").Value = 200
$ws.Range("Q451").Value = 1
$ws.Range("R451").Value = "Hortaliza"
